$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (row 2 is the first data row, so this becomes
# the new "latest" weekly data point), which shifts all existing data rows
# (previously 3..28) down to (4..29).
$ws.Rows("3").Insert()

# Populate the newly inserted row 3 with this week's data.
$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(3, 3).Value = "Los Lagos"
$ws.Cells.Item(3, 4).Value = 45149
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 6).Value = 100112035
$ws.Cells.Item(3, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 120
$ws.Cells.Item(3, 11).Value = 25000
$ws.Cells.Item(3, 12).Value = 25000
$ws.Cells.Item(3, 13).Value = 25000
$ws.Cells.Item(3, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(3, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(3, 16).Value = 1667
$ws.Cells.Item(3, 17).Value = 15
$ws.Cells.Item(3, 18).Value = "Hortaliza"
